$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 7011
$ws.Range("J3").Value = 7400
$ws.Range("J4").Value = 1615
$ws.Range("J6").Value = 10012
$ws.Range("J7").Value = 26616

$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("J5").Value = 4
$ws.Range("J6").Value = 12

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J3").Value = 495
$ws.Range("J4").Value = 89
$ws.Range("J7").Value = 1669

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J6").Value = 143
$ws.Range("J7").Value = 533

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 276
$ws.Range("J3").Value = 397
$ws.Range("J4").Value = 55
$ws.Range("J6").Value = 424
$ws.Range("J7").Value = 1202

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J3").Value = 138
$ws.Range("J7").Value = 387

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 244
$ws.Range("J6").Value = 238

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J2").Value = 111
$ws.Range("J7").Value = 407

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J5").Value = 79
$ws.Range("J6").Value = 202
$ws.Range("J8").Value = 1669
$ws.Range("J9").Value = 138
$ws.Range("J10").Value = 197
$ws.Range("J21").Value = 73
$ws.Range("J23").Value = 244
$ws.Range("J24").Value = 85
$ws.Range("J27").Value = 160
$ws.Range("J29").Value = 1427
$ws.Range("J32").Value = 43
$ws.Range("J33").Value = 1202
$ws.Range("J34").Value = 121
$ws.Range("J38").Value = 12
$ws.Range("J40").Value = 59
$ws.Range("J42").Value = 1149
$ws.Range("J48").Value = 301
$ws.Range("J51").Value = 329
$ws.Range("J54").Value = 517
$ws.Range("J63").Value = 83
$ws.Range("J64").Value = 176
$ws.Range("J67").Value = 992
$ws.Range("J73").Value = 255
$ws.Range("J79").Value = 740
$ws.Range("J83").Value = 533
$ws.Range("J85").Value = 1096
$ws.Range("J88").Value = 285
$ws.Range("J91").Value = 307
$ws.Range("J93").Value = 111
$ws.Range("J95").Value = 387
$ws.Range("J96").Value = 290
$ws.Range("J99").Value = 407
$ws.Range("J101").Value = 26616

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 254
$ws.Range("J3").Value = 369
$ws.Range("J4").Value = 67
$ws.Range("J6").Value = 275
$ws.Range("J7").Value = 992

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J2").Value = 128
$ws.Range("J3").Value = 105
$ws.Range("J4").Value = 39
$ws.Range("J7").Value = 517

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J4").Value = 77
$ws.Range("J7").Value = 1427

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J2").Value = 50
$ws.Range("J7").Value = 301

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("J2").Value = 61
$ws.Range("J6").Value = 73
$ws.Range("J7").Value = 202

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("J3").Value = 28
$ws.Range("J6").Value = 118

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J3").Value = 229
$ws.Range("J6").Value = 612
$ws.Range("J7").Value = 1149

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J6").Value = 112
$ws.Range("J7").Value = 197

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J2").Value = 30
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 85

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J6").Value = 66
$ws.Range("J7").Value = 244

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J2").Value = 88
$ws.Range("J7").Value = 290

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J3").Value = 127
$ws.Range("J7").Value = 307

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("J2").Value = 12
$ws.Range("J3").Value = 13
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 208
$ws.Range("J7").Value = 740

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J6").Value = 62
$ws.Range("J7").Value = 176

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("J2").Value = 31
$ws.Range("J7").Value = 111

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 121

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J6").Value = 47
$ws.Range("J7").Value = 138

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J2").Value = 82
$ws.Range("J7").Value = 255

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 148
$ws.Range("J7").Value = 285

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("J2").Value = 13
$ws.Range("J7").Value = 43

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("J6").Value = 39
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J3").Value = 39
$ws.Range("J4").Value = 20
$ws.Range("J7").Value = 160

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J3").Value = 87
$ws.Range("J4").Value = 28
$ws.Range("J6").Value = 134
$ws.Range("J7").Value = 329

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 293
$ws.Range("J3").Value = 392
$ws.Range("J4").Value = 71
$ws.Range("J6").Value = 314
$ws.Range("J7").Value = 1096

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J2").Value = 27
$ws.Range("J4").Value = 8

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("J4").Value = 7
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 59
